$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells: role, manager (C1, D1) ---
$ws.Range("C1").Value = "role"
$ws.Range("D1").Value = "manager"

# --- New user names first (column A), rows 4 and 5 ---
$ws.Range("A4").Value = "gamma"
$ws.Range("A5").Value = "delta"

# --- role column (C2:C5) ---
$ws.Range("C2").Value = "manager"
$ws.Range("C3").Value = "employee"
$ws.Range("C4").Value = "admin"
$ws.Range("C5").Value = "superuser"

# --- tenure header (E1) ---
$ws.Range("E1").Value = "tenure"

# --- remaining gender values for the new users ---
$ws.Range("B4").Value = "male"
$ws.Range("B5").Value = "female"

# --- manager column (D3:D5) ---
$ws.Range("D3").Value = "alpha"
$ws.Range("D4").Value = "alpha"
$ws.Range("D5").Value = "alpha"

# --- tenure values (E2:E5) ---
$ws.Range("E2").Value = 103
$ws.Range("E2").NumberFormat = "General"
$ws.Range("E3").Value = 30
$ws.Range("E3").NumberFormat = "General"
$ws.Range("E4").Value = 12
$ws.Range("E4").NumberFormat = "General"
$ws.Range("E5").Value = 94
$ws.Range("E5").NumberFormat = "General"

# --- Column E width (Excel auto "best fit" width after entering the tenure numbers) ---
$ws.Range("E:E").ColumnWidth = 10.140625

# --- Page setup (paper size / orientation) as added by Excel on save ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Update selection to match the final active cell in the authored file ---
$ws.Range("E5").Select()
